$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two score values; the SUM formula in C51 will recalc automatically.
$ws.Range("C8").Value = 15
$ws.Range("C9").Value = 27

# Move the selection/active cell to C8 (also clears the scrolled topLeftCell position).
$ws.Range("C8").Select()
